$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("I2").Value = 0.4510973247164258
$ws.Range("J2").Value = 0.4510973247164258
$ws.Range("K2").Value = 1
$ws.Range("L2").Value = 0.3333333333333333
$ws.Range("M2").Value = 0.190825
$ws.Range("N2").Value = 0.572475
$ws.Range("O2").Value = 0.6793270274792366
$ws.Range("P2").Value = 0.6793270274792366
$ws.Range("Q2").Value = 0.03880458179166667
$ws.Range("R2").Value = 0.349241236125
$ws.Range("S2").Value = 0.3064426047034455
$ws.Range("T2").Value = 0.3064426047034455

# Row 3
$ws.Range("I3").Value = 0.4510973247164258
$ws.Range("J3").Value = 0.4510973247164258
$ws.Range("M3").Value = 0.09007799999999999
$ws.Range("O3").Value = 0.3206729725207633
$ws.Range("P3").Value = 0.3206729725207634
$ws.Range("S3").Value = 0.1446547200129803
$ws.Range("T3").Value = 0.1446547200129803

# Row 4
$ws.Range("E4").Value = 2
$ws.Range("F4").Value = 0.6666666666666666
$ws.Range("G4").Value = 0.2474416666666667
$ws.Range("H4").Value = 0.742325
$ws.Range("I4").Value = 0.5489026752835741
$ws.Range("J4").Value = 0.5489026752835741
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.190825
$ws.Range("N4").Value = 0.572475
$ws.Range("O4").Value = 0.6793270274792366
$ws.Range("P4").Value = 0.6793270274792366
$ws.Range("Q4").Value = 0.04721805604166666
$ws.Range("R4").Value = 0.424962504375
$ws.Range("S4").Value = 0.372884422775791
$ws.Range("T4").Value = 0.372884422775791

# Row 5
$ws.Range("E5").Value = 2
$ws.Range("F5").Value = 0.6666666666666666
$ws.Range("G5").Value = 0.2474416666666667
$ws.Range("H5").Value = 0.742325
$ws.Range("I5").Value = 0.5489026752835741
$ws.Range("J5").Value = 0.5489026752835741
$ws.Range("M5").Value = 0.09007799999999999
$ws.Range("O5").Value = 0.3206729725207633
$ws.Range("P5").Value = 0.3206729725207634
$ws.Range("Q5").Value = 0.02228905045
$ws.Range("R5").Value = 0.20060145405
$ws.Range("S5").Value = 0.176018252507783
$ws.Range("T5").Value = 0.1760182525077831
